# This script applies a cyclic rotation of the "species observation" data
# (columns A, B, D, E, F, G, H, Q, R, AC) across rows 12-16:
#   row 12 data -> row 16
#   row 13 data -> row 14
#   row 14 data -> row 12
#   row 15 data -> row 13
#   row 16 data -> row 15
# All other columns for these rows stay as-is (they already hold identical
# values across the block), so only the columns above need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","Q","R","AC")

# Capture the current ("before") values for rows 12-16 for the columns that
# move, using the Value() accessor (needed to actually invoke the getter
# in this host rather than stringify the property descriptor).
$before = @{}
foreach ($r in 12..16) {
    $rowData = @{}
    foreach ($col in $cols) {
        $addr = "${col}${r}"
        $rowData[$col] = $ws.Range($addr).Value()
    }
    $before[$r] = $rowData
}

# Mapping of source row (old position) -> destination row (new position).
$mapping = @{
    12 = 16
    13 = 14
    14 = 12
    15 = 13
    16 = 15
}

foreach ($src in $mapping.Keys) {
    $dst = $mapping[$src]
    $srcData = $before[$src]
    foreach ($col in $cols) {
        $addr = "${col}${dst}"
        $ws.Range($addr).Value = $srcData[$col]
    }
}
